# Fix some dirty files:
# 1. Correct a truncated shared string in the 保險 (Insurance) sheet.
# 2. Correct the running "index" numbering in several sheets (股票, 基金受益憑證, 保險, 債務)
#    which had been accidentally bumped too high.

$wb = $excel.ActiveWorkbook

# --- 1. Fix insurance product name text ---
$wsInsurance = $wb.Worksheets.Item("保險")
$oldText = "全球人壽全球104终身壽險"
$newText = "全球人壽全球104终身壽險甲型"

$found = $wsInsurance.Cells.Find($oldText)
if ($found -ne $null) {
    $found.Value = $newText
}

# --- 2. Fix "index" numbering ---

# 股票 (Stock) sheet: rows 2-4, columns A and N
$wsStock = $wb.Worksheets.Item("股票")
$wsStock.Range("A2").Value = 76
$wsStock.Range("N2").Value = 76
$wsStock.Range("A3").Value = 77
$wsStock.Range("N3").Value = 77
$wsStock.Range("A4").Value = 78
$wsStock.Range("N4").Value = 78

# 基金受益憑證 (Fund) sheet: rows 2-3, columns A and O
$wsFund = $wb.Worksheets.Item("基金受益憑證")
$wsFund.Range("A2").Value = 88
$wsFund.Range("O2").Value = 88
$wsFund.Range("A3").Value = 89
$wsFund.Range("O3").Value = 89

# 保險 (Insurance) sheet: rows 2-9, columns A and K
$wsInsurance.Range("A2").Value = 103
$wsInsurance.Range("K2").Value = 103
$wsInsurance.Range("A3").Value = 104
$wsInsurance.Range("K3").Value = 104
$wsInsurance.Range("A4").Value = 105
$wsInsurance.Range("K4").Value = 105
$wsInsurance.Range("A5").Value = 106
$wsInsurance.Range("K5").Value = 106
$wsInsurance.Range("A6").Value = 107
$wsInsurance.Range("K6").Value = 107
$wsInsurance.Range("A7").Value = 108
$wsInsurance.Range("K7").Value = 108
$wsInsurance.Range("A8").Value = 109
$wsInsurance.Range("K8").Value = 109
$wsInsurance.Range("A9").Value = 110
$wsInsurance.Range("K9").Value = 110

# 債務 (Debt) sheet: rows 2-6, columns A and N
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Range("A2").Value = 120
$wsDebt.Range("N2").Value = 120
$wsDebt.Range("A3").Value = 121
$wsDebt.Range("N3").Value = 121
$wsDebt.Range("A4").Value = 122
$wsDebt.Range("N4").Value = 122
$wsDebt.Range("A5").Value = 123
$wsDebt.Range("N5").Value = 123
$wsDebt.Range("A6").Value = 124
$wsDebt.Range("N6").Value = 124
